$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 5: "Déploiement 1.4.1446.0 + RCH" ---
$ws.Range("B5").Value = "11/14/2014"
$ws.Range("C5").Value = "Déploiement 1.4.1446.0 + RCH"

# --- Row 25: "Chasse au bug, analyse registre" ---
$ws.Range("B25").Value = "11/11/2014"
$ws.Range("C25").Value = "Chasse au bug, analyse registre"

# --- Row 26: "Gestion des déf. De groupes ""fonctions""" ---
$ws.Range("B26").Value = "11/13/2014"
$ws.Range("C26").Value = 'Gestion des déf. De groupes "fonctions"'

# --- Row 27: "Analyse Registres" ---
$ws.Range("B27").Value = "11/19/2014"
$ws.Range("C27").Value = "Analyse Registres"

# --- Row 6: "Export BN + maintenance" ---
$ws.Range("B6").Value = "11/20/2014"
$ws.Range("C6").Value = "Export BN + maintenance"

# --- Row 28: "Analyse, médlisation, modifications" ---
$ws.Range("B28").Value = "11/25/2014"
$ws.Range("C28").Value = "Analyse, médlisation, modifications"

# --- Row 29: "Collaborateurs AIDER" ---
$ws.Range("B29").Value = "11/26/2014"
$ws.Range("C29").Value = "Collaborateurs AIDER"

# --- Row 30: "Utilisateurs AIDER" ---
$ws.Range("B30").Value = "11/27/2014"
$ws.Range("C30").Value = "Utilisateurs AIDER"

# --- Row 7: "Divers, réponses aux mails, etc." (reuses an existing shared string) ---
$ws.Range("B7").Value = "11/24/2014"
$ws.Range("C7").Value = "Divers, réponses aux mails, etc."

# --- Numeric (time) values ---
$ws.Range("E5").Value = 0.055555555555555552
$ws.Range("E5").NumberFormat = $ws.Range("D5").NumberFormat

$ws.Range("E6").Value = 0.041666666666666664
$ws.Range("E6").NumberFormat = $ws.Range("D6").NumberFormat

$ws.Range("D7").Value = 0.0069444444444444441

$ws.Range("D25").Value = 0.16666666666666666
$ws.Range("D26").Value = 0.16666666666666666
$ws.Range("D27").Value = 0.10416666666666667
$ws.Range("D28").Value = 0.20833333333333334
$ws.Range("D29").Value = 0.29166666666666669
$ws.Range("D30").Value = 0.16666666666666666

# --- Update the selection to match the author's final cursor position ---
$ws.Range("B31").Select()
